$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $rng = $ws.Range($cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue "D2" "44.320.36"
Set-TextValue "E2" "  +2.58%  "
Set-TextValue "D3" "2.367.17"
Set-TextValue "E3" "  +0.47%  "
Set-TextValue "E4" "  +0.01%  "
Set-TextValue "E5" "  +3.90%  "
Set-TextValue "D6" "238.76"
Set-TextValue "E6" "  +2.69%  "
Set-TextValue "D7" "73.44"
Set-TextValue "E7" "  +7.94%  "
Set-TextValue "E8" "  -0.01%  "
Set-TextValue "D9" "0.551"
Set-TextValue "E9" "  +20.24%  "
Set-TextValue "E10" "  +6.75%  "
Set-TextValue "D11" "30.21"
Set-TextValue "E11" "  +14.53%  "
Set-TextValue "E12" "  +2.15%  "
Set-TextValue "D13" "2.716.44"
Set-TextValue "E13" "  +0.34%  "
Set-TextValue "D14" "16.86"
Set-TextValue "E14" "  +7.80%  "
Set-TextValue "D15" "6.82"
Set-TextValue "E15" "  +9.34%  "
Set-TextValue "E16" "  +7.84%  "
Set-TextValue "D17" "2.370.65"
Set-TextValue "E17" "  +0.46%  "
Set-TextValue "D18" "44.483.04"
Set-TextValue "E18" "  +2.92%  "
Set-TextValue "E19" "  +4.70%  "
Set-TextValue "D20" "77.38"
Set-TextValue "E20" "  +4.56%  "
Set-TextValue "D21" "6.48"
Set-TextValue "E21" "  +3.90%  "
Set-TextValue "D22" "254.69"
Set-TextValue "E22" "  +2.47%  "
Set-TextValue "D23" "3.87"
Set-TextValue "E23" "  -4.36%  "
Set-TextValue "E24" "  +0.01%  "
Set-TextValue "E25" "  +2.66%  "
Set-TextValue "D26" "10.40"
Set-TextValue "E26" "  +4.92%  "
Set-TextValue "E27" "  -1.28%  "
Set-TextValue "D28" "22.66"
Set-TextValue "E28" "  +1.45%  "
Set-TextValue "D29" "1.60"
Set-TextValue "E29" "  +4.60%  "
Set-TextValue "D30" "174.00"
Set-TextValue "E30" "  +0.37%  "
Set-TextValue "E31" "  +2.65%  "
Set-TextValue "E32" "  +5.24%  "
Set-TextValue "D33" "0.0742"
Set-TextValue "E33" "  +6.89%  "
Set-TextValue "D34" "5.21"
Set-TextValue "E34" "  +3.95%  "
Set-TextValue "D35" "5.23"
Set-TextValue "E35" "  +3.00%  "
Set-TextValue "E36" "  +7.25%  "
Set-TextValue "E37" "  -2.91%  "
Set-TextValue "D38" "6.47"
Set-TextValue "E38" "  -0.36%  "
Set-TextValue "D39" "0.0272"
Set-TextValue "E39" "  +6.52%  "
Set-TextValue "D40" "19.42"
Set-TextValue "E40" "  +7.07%  "
Set-TextValue "E41" "  +0.13%  "
Set-TextValue "D42" "8.85"
Set-TextValue "E42" "  -0.94%  "
Set-TextValue "E43" "  +3.25%  "
Set-TextValue "D44" "0.0989"
Set-TextValue "E44" "  +4.16%  "
Set-TextValue "E45" "  +1.54%  "
Set-TextValue "D46" "0.186"
Set-TextValue "E46" "  +12.63%  "
Set-TextValue "B47" "Aave"
Set-TextValue "C47" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D47" "99.09"
Set-TextValue "E47" "  +0.49%  "
Set-TextValue "B48" "FTXToken"
Set-TextValue "C48" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D48" "4.48"
Set-TextValue "E48" "  +0.89%  "
Set-TextValue "D49" "2.38"
Set-TextValue "E49" "  +4.53%  "
Set-TextValue "D50" "1.445.03"
Set-TextValue "E50" "  +0.07%  "
Set-TextValue "D51" "2.590.61"
Set-TextValue "E51" "  +0.41%  "
